# Add a new "Экзамен" (Exam) score column.
# The header goes into I2, and most rows get their score in column I,
# except row 10 whose score ended up in column J.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Cells.Item(2, 9).Value = "Экзамен"

# Exam scores for each student row (row 21 - "отчислен" / expelled - has none)
$ws.Cells.Item(4, 9).Value  = 5
$ws.Cells.Item(5, 9).Value  = 4
$ws.Cells.Item(6, 9).Value  = 5
$ws.Cells.Item(7, 9).Value  = 5
$ws.Cells.Item(8, 9).Value  = 3
$ws.Cells.Item(9, 9).Value  = 4
$ws.Cells.Item(10, 10).Value = 5   # row 10 uses column J instead of I
$ws.Cells.Item(11, 9).Value = 5
$ws.Cells.Item(12, 9).Value = 5
$ws.Cells.Item(13, 9).Value = 3
$ws.Cells.Item(14, 9).Value = 4
$ws.Cells.Item(15, 9).Value = 5
$ws.Cells.Item(16, 9).Value = 4
$ws.Cells.Item(17, 9).Value = 5
$ws.Cells.Item(18, 9).Value = 4
$ws.Cells.Item(19, 9).Value = 5
$ws.Cells.Item(20, 9).Value = 4
$ws.Cells.Item(22, 9).Value = 5
$ws.Cells.Item(23, 9).Value = 4
$ws.Cells.Item(24, 9).Value = 4
$ws.Cells.Item(25, 9).Value = 5

# Update the active selection to match the target state
$ws.Range("I12").Select()
